# Updates the crypto price/volume table to the latest scraped values.
# (Generated from the authoritative cell-level diff.)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value as plain TEXT, without leaving any lasting
# NumberFormat/Style change on the cell (price strings like "1.002" or
# "215.03" would otherwise be auto-coerced by Excel into numbers,
# truncating meaningful trailing zeros / losing the "t=inlineStr" shape).
function Set-TextValue($addr, $value) {
    $r = $ws.Range($addr)
    $r.Style = "Normal"
    $r.NumberFormat = "@"
    $r.Value = $value
    $r.Style = "Normal"
}


# Row 2
Set-TextValue 'D2' '25.883.30'
Set-TextValue 'E2' '  -1.01%  '

# Row 3
Set-TextValue 'D3' '1.637.06'
Set-TextValue 'E3' '  -0.88%  '

# Row 4
Set-TextValue 'E4' '  -0.34%  '

# Row 5
Set-TextValue 'D5' '215.03'
Set-TextValue 'E5' '  -0.18%  '

# Row 6
Set-TextValue 'D6' '0.5033'
Set-TextValue 'E6' '  -1.83%  '

# Row 7
Set-TextValue 'D7' '1.002'
Set-TextValue 'E7' '  -0.34%  '

# Row 8
Set-TextValue 'D8' '0.2563'
Set-TextValue 'E8' '  -1.14%  '

# Row 9
Set-TextValue 'D9' '0.06387'
Set-TextValue 'E9' '  -0.82%  '

# Row 10
Set-TextValue 'D10' '19.66'
Set-TextValue 'E10' '  -1.32%  '

# Row 11
Set-TextValue 'D11' '0.07729'
Set-TextValue 'E11' '  -0.66%  '

# Row 12
Set-TextValue 'D12' '4.256'
Set-TextValue 'E12' '  -0.80%  '

# Row 13
$ws.Range('B13').Value = 'WrappedEther'
$ws.Range('C13').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
Set-TextValue 'D13' '1.637.89'
Set-TextValue 'E13' '  -0.95%  '

# Row 14
$ws.Range('B14').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C14').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
Set-TextValue 'D14' '1.863.57'
Set-TextValue 'E14' '  -0.80%  '

# Row 15
Set-TextValue 'D15' '0.5440'
Set-TextValue 'E15' '  -1.38%  '

# Row 16
Set-TextValue 'D16' '0.0₅7895'
Set-TextValue 'E16' '  -1.31%  '

# Row 17
Set-TextValue 'D17' '64.24'
Set-TextValue 'E17' '  +0.25%  '

# Row 18
Set-TextValue 'D18' '25.892.21'

# Row 19
Set-TextValue 'D19' '1.002'
Set-TextValue 'E19' '  -0.34%  '

# Row 20
Set-TextValue 'D20' '203.33'
Set-TextValue 'E20' '  -3.46%  '

# Row 21
Set-TextValue 'D21' '4.374'
Set-TextValue 'E21' '  -0.54%  '

# Row 22
Set-TextValue 'D22' '9.898'
Set-TextValue 'E22' '  -1.57%  '

# Row 23
Set-TextValue 'D23' '5.973'
Set-TextValue 'E23' '  -1.16%  '

# Row 24
Set-TextValue 'D24' '1.003'
Set-TextValue 'E24' '  -0.27%  '

# Row 25
Set-TextValue 'D25' '1.934'
Set-TextValue 'E25' '  +10.57%  '

# Row 26
Set-TextValue 'D26' '140.97'
Set-TextValue 'E26' '  -1.72%  '

# Row 27
Set-TextValue 'D27' '0.1137'
Set-TextValue 'E27' '  -3.26%  '

# Row 28
Set-TextValue 'D28' '15.69'
Set-TextValue 'E28' '  -0.73%  '

# Row 29
Set-TextValue 'D29' '6.747'
Set-TextValue 'E29' '  -3.15%  '

# Row 30
Set-TextValue 'D30' '1.242'
Set-TextValue 'E30' '  +0.15%  '

# Row 31
Set-TextValue 'D31' '0.04944'
Set-TextValue 'E31' '  -2.92%  '

# Row 32
Set-TextValue 'D32' '3.272'
Set-TextValue 'E32' '  -2.43%  '

# Row 33
Set-TextValue 'D33' '3.182'
Set-TextValue 'E33' '  -1.14%  '

# Row 34
Set-TextValue 'D34' '1.543'
Set-TextValue 'E34' '  -1.06%  '

# Row 35
Set-TextValue 'D35' '2.370'
Set-TextValue 'E35' '  +0.84%  '

# Row 36
Set-TextValue 'E36' '  -4.04%  '

# Row 37
Set-TextValue 'D37' '0.8905'
Set-TextValue 'E37' '  -3.52%  '

# Row 38
Set-TextValue 'D38' '1.157.96'
Set-TextValue 'E38' '  -0.09%  '

# Row 39
Set-TextValue 'E39' '  -1.92%  '

# Row 40
Set-TextValue 'D40' '0.01564'
Set-TextValue 'E40' '  -1.36%  '

# Row 41
Set-TextValue 'D41' '1.002'
Set-TextValue 'E41' '  -0.33%  '

# Row 42
Set-TextValue 'D42' '5.636'
Set-TextValue 'E42' '  -0.29%  '

# Row 43
$ws.Range('B43').Value = 'Quant'
$ws.Range('C43').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
Set-TextValue 'D43' '99.84'
Set-TextValue 'E43' '  -0.24%  '

# Row 44
$ws.Range('B44').Value = 'TrustWalletToken'
$ws.Range('C44').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
Set-TextValue 'D44' '0.8060'
Set-TextValue 'E44' '  -2.05%  '

# Row 45
Set-TextValue 'D45' '1.775.26'
Set-TextValue 'E45' '  -0.79%  '

# Row 46
Set-TextValue 'E46' '  +0.67%  '

# Row 47
Set-TextValue 'D47' '0.4530'
Set-TextValue 'E47' '  -0.45%  '

# Row 48
$ws.Range('B48').Value = 'Frax'
$ws.Range('C48').Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
Set-TextValue 'D48' '1.001'
Set-TextValue 'E48' '  -0.39%  '

# Row 49
$ws.Range('B49').Value = 'Aave'
$ws.Range('C49').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
Set-TextValue 'D49' '54.88'
Set-TextValue 'E49' '  -1.05%  '

# Row 50
Set-TextValue 'E50' '  -0.60%  '

# Row 51
Set-TextValue 'D51' '1.003'
Set-TextValue 'E51' '  -0.29%  '

Write-Output "cryptos list updated"
